# Edit described by the commit "change the most problematic sentences".
#
# The workbook has three sheets:
#   1 = lemmas
#   2 = sentence_group_1  (subject "Student" for the zapsat/row16 lemma)
#   3 = sentence_group_2  (subject "Vratny" for the zapsat/row16 lemma)
#
# For the lemma "zapsat" (row 16) the example sentences that used
# "do sesitu"/"do kalendare" are replaced with new sentences using
# "pri odchodu"/"u vchodu". The order below matches the order the
# shared strings were written in the original edit.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("E16").Value = "Vrátný zapsal návštěvu u vchodu"
$ws3.Range("F16").Value = "Vrátný se zapsal návštěvu u vchodu"
$ws3.Range("C16").Value = "Vrátný se zapsal u vchodu"
$ws3.Range("D16").Value = "Vrátný zapsal u vchodu"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E16").Value = "Student zapsal výsledky při odchodu"
$ws2.Range("F16").Value = "Student se zapsal výsledky při odchodu"
$ws2.Range("C16").Value = "Student se zapsal při odchodu"
$ws2.Range("D16").Value = "Student zapsal při odchodu"

# Highlight (clear to white) column D - "refl-ngr" - for the rows that were
# flagged as problematic.
$ws2.Range("D10").Interior.Color = 16777215

$ws3.Range("D10").Interior.Color = 16777215
$ws3.Range("D11").Interior.Color = 16777215
$ws3.Range("D12").Interior.Color = 16777215
$ws3.Range("D13").Interior.Color = 16777215
$ws3.Range("D14").Interior.Color = 16777215
$ws3.Range("D15").Interior.Color = 16777215
$ws3.Range("D16").Interior.Color = 16777215

# Update sheet selections/active tab: sentence_group_2 was active before and
# becomes inactive (selection left at C25); sentence_group_1 becomes the
# active tab with selection at E22.
$ws3.Select()
$ws3.Range("C25").Select()

$ws2.Select()
$ws2.Range("E22").Select()
